# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures scraped
# from coinranking.com, and fixes the OKB / VeChain row ordering (rows 43-44).
#
# NOTE: values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr cells) instead of coercing
# strings like "0.750" or "69.70" into numbers and dropping trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-Text "D2" "58.695.46"
Set-Text "E2" "  +1.92%  "

# Row 3 - Ethereum
Set-Text "D3" "3.164.26"
Set-Text "E3" "  +1.59%  "

# Row 4 - TetherUSD
Set-Text "E4" "  -0.02%  "

# Row 5 - BNB
Set-Text "D5" "529.43"
Set-Text "E5" "  -0.19%  "

# Row 6 - Solana
Set-Text "D6" "139.85"

# Row 7 - USDC
Set-Text "E7" "  -0.08%  "

# Row 8 - XRP
Set-Text "E8" "  +14.23%  "

# Row 9 - Toncoin
Set-Text "E9" "  -0.05%  "

# Row 10 - Cardano
Set-Text "D10" "0.439"
Set-Text "E10" "  +6.74%  "

# Row 11 - Dogecoin
Set-Text "E11" "  +3.98%  "

# Row 12 - TRON
Set-Text "E12" "  +2.51%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Text "D13" "3.709.39"

# Row 14 - Avalanche
Set-Text "D14" "25.78"
Set-Text "E14" "  +0.66%  "

# Row 15 - ShibaInu
Set-Text "E15" "  +3.77%  "

# Row 16 - WrappedBTC
Set-Text "D16" "58.738.85"
Set-Text "E16" "  +1.81%  "

# Row 17 - WrappedEther
Set-Text "D17" "3.177.22"
Set-Text "E17" "  +1.97%  "

# Row 18 - Polkadot
Set-Text "D18" "6.25"
Set-Text "E18" "  +3.55%  "

# Row 19 - Chainlink
Set-Text "D19" "12.99"
Set-Text "E19" "  +2.45%  "

# Row 20 - BitcoinCash
Set-Text "D20" "376.07"
Set-Text "E20" "  +4.24%  "

# Row 21 - Uniswap
Set-Text "E21" "  +0.17%  "

# Row 22 - Dai
Set-Text "E22" "  +0.24%  "

# Row 23 - Polygon
Set-Text "E23" "  +4.95%  "

# Row 24 - Litecoin
Set-Text "D24" "69.70"
Set-Text "E24" "  +0.97%  "

# Row 25 - Kaspa
Set-Text "D25" "0.168"
Set-Text "E25" "  +0.56%  "

# Row 26 - Binance-PegBSC-USD
Set-Text "E26" "  +0.02%  "

# Row 27 - InternetComputer(DFINITY)
Set-Text "D27" "8.32"
Set-Text "E27" "  +14.18%  "

# Row 28 - PEPE
Set-Text "E28" "  +0.02%  "

# Row 29 - EthereumClassic
Set-Text "D29" "22.43"
Set-Text "E29" "  +5.00%  "

# Row 30 - PancakeSwap
Set-Text "E30" "  +1.02%  "

# Row 31 - RenderToken
Set-Text "E31" "  -0.83%  "

# Row 32 - NEARProtocol
Set-Text "D32" "5.15"
Set-Text "E32" "  +0.78%  "

# Row 33 - Fetch.AI
Set-Text "E33" "  +1.25%  "

# Row 34 - Aptos
Set-Text "E34" "  +4.39%  "

# Row 35 - Monero
Set-Text "D35" "156.87"
Set-Text "E35" "  -1.44%  "

# Row 36 - ImmutableX
Set-Text "E36" "  +4.22%  "

# Row 37 - EnergySwap
Set-Text "D37" "25.04"
Set-Text "E37" "  -1.60%  "

# Row 38 - Maker
Set-Text "D38" "2.689.27"
Set-Text "E38" "  +8.13%  "

# Row 39 - Hedera
Set-Text "E39" "  +3.95%  "

# Row 40 - Stacks
Set-Text "E40" "  +1.48%  "

# Row 41 - Filecoin
Set-Text "E41" "  +6.92%  "

# Row 42 - Mantle
Set-Text "D42" "0.724"
Set-Text "E42" "  +4.05%  "

# Row 43 - was OKB, now VeChain (rows 43/44 swap coin ordering)
Set-Text "B43" "VeChain"
Set-Text "C43" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Text "D43" "0.0290"
Set-Text "E43" "  +7.51%  "

# Row 44 - was VeChain, now OKB
Set-Text "B44" "OKB"
Set-Text "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-Text "D44" "39.14"
Set-Text "E44" "  +3.76%  "

# Row 45 - FirstDigitalUSD
Set-Text "E45" "  +0.00%  "

# Row 46 - RenzoRestakedETH
Set-Text "D46" "3.207.17"
Set-Text "E46" "  +1.61%  "

# Row 47 - Stellar
Set-Text "D47" "0.103"
Set-Text "E47" "  +13.52%  "

# Row 48 - Cosmos
Set-Text "D48" "6.21"
Set-Text "E48" "  +2.28%  "

# Row 49 - ONDO
Set-Text "D49" "0.978"
Set-Text "E49" "  +0.08%  "

# Row 50 - InjectiveProtocol
Set-Text "D50" "20.07"
Set-Text "E50" "  +2.08%  "

# Row 51 - SuiNetwork
Set-Text "D51" "0.750"
Set-Text "E51" "  +1.41%  "
